$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "zou123"
$ws.Range("C6").Value = "zou123456"

$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("C6").HorizontalAlignment = -4108
$ws.Range("C6").VerticalAlignment = -4108

$ws.Range("C7").Select()
